$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 935.0898621313619
$ws.Range("D2").Value = 4.675449310656809
$ws.Range("E2").Value = 6.493679598134458
$ws.Range("B3").Value = 701.3173965985216
$ws.Range("D3").Value = 3.5065869829926077
$ws.Range("E3").Value = 4.870259698600844
$ws.Range("B4").Value = 124.6786482841816
$ws.Range("D4").Value = 0.623393241420908
$ws.Range("E4").Value = 0.8658239464179278
$ws.Range("B5").Value = 1761.0859070140648
$ws.Range("D5").Value = 8.805429535070324
$ws.Range("E5").Value = 12.22976324315323
$ws.Range("B7").Value = 258.0848019482559
$ws.Range("D7").Value = 1.2904240097412796
$ws.Range("E7").Value = 1.7922555690851105
$ws.Range("B8").Value = 43.01413365804265
$ws.Range("D8").Value = 0.21507066829021326
$ws.Range("E8").Value = 0.2987092615141851
$ws.Range("B9").Value = 301.0989356062986
$ws.Range("D9").Value = 1.505494678031493
$ws.Range("E9").Value = 2.0909648305992956
$ws.Range("B11").Value = 276.30763791891854
$ws.Range("C11").Value = 128.47295268831033
$ws.Range("D11").Value = 1.3815381895945926
$ws.Range("E11").Value = 1.918803041103601
$ws.Range("B13").Value = 171.2523556921675
$ws.Range("C13").Value = 79.62608618534725
$ws.Range("D13").Value = 0.8562617784608375
$ws.Range("E13").Value = 1.1892524700844964
$ws.Range("B18").Value = 108.74524586452637
$ws.Range("C18").Value = 50.5625647276955
$ws.Range("D18").Value = 0.5437262293226318
$ws.Range("E18").Value = 0.7551753185036554
$ws.Range("B20").Value = 395.5267444413777
$ws.Range("C20").Value = 183.90547980613496
$ws.Range("D20").Value = 1.9776337222068885
$ws.Range("E20").Value = 2.7467135030651226
$ws.Range("B21").Value = 282.65787645581497
$ws.Range("C21").Value = 131.42558150905103
$ws.Range("D21").Value = 1.4132893822790749
$ws.Range("E21").Value = 1.9629019198320483
$ws.Range("B22").Value = 760.1986646632856
$ws.Range("C22").Value = 353.4645940828549
$ws.Range("D22").Value = 3.800993323316428
$ws.Range("E22").Value = 5.279157393495039
$ws.Range("B25").Value = 3207.4363910670936
$ws.Range("C25").Value = 1491.3406912090054
$ws.Range("D25").Value = 16.03718195533547
$ws.Range("E25").Value = 22.27386382685482
$ws.Range("B27").Value = 1446.350484053029
$ws.Range("C27").Value = 672.5001114988607
$ws.Range("D27").Value = 7.2317524202651455
$ws.Range("E27").Value = 10.04410058370159
